# Update the "working days" counters on the Config sheet (October run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: bump the numeric counter 42 -> 43
$ws.Range("B2").Value = "43"

# D2: bump the automation tag Automation3 -> Automation2
$ws.Range("D2").Value = "Automation2"

# Update the active selection to match the saved workbook state (E8)
$ws.Range("E8").Select()
